$d = $word.ActiveDocument

# 1. Remove the prototype image (run) at the top of the document, keeping
#    the (now empty) paragraph, then drop a fresh "_GoBack" bookmark in it.
$shape = $d.InlineShapes.Item(1)
$shape.Delete()

$p1 = $d.Paragraphs.Item(1)
$p1.Range.Bookmarks.Add("_GoBack", $p1.Range)

# 2. Remove the stray <w:lastRenderedPageBreak/> before "TESTE No."
$d.Content.Find.Execute("TESTE No.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "TESTE No.", 2)

# 3. Update the sectPr header/footer relationship ids (rId9/rId10 -> rId8/rId9)
